$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (ID "H 72"), shifting all rows below up by one.
$ws.Rows.Item(2).Delete()
